$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 4. This shifts the existing row 4 (and all
# rows below it) down by one row, preserving all of their values and
# formatting -- matching the "row 5 becomes old row 4", ... , "row 66
# becomes old row 65" shift seen in the target diff.
$ws.Rows.Item(4).Insert()

# Populate the newly inserted row 4 with this week's new record (a fresh
# price/date entry, otherwise identical to the rest of the series).
$ws.Cells.Item(4,1).Value = 4
$ws.Cells.Item(4,2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(4,3).Value = "Los Lagos"
$ws.Cells.Item(4,4).Value = 45083
$ws.Cells.Item(4,5).Value = 10
$ws.Cells.Item(4,6).Value = 100112012
$ws.Cells.Item(4,7).Value = "Espinaca"
$ws.Cells.Item(4,8).Value = "Sin especificar"
$ws.Cells.Item(4,9).Value = "Primera"
$ws.Cells.Item(4,10).Value = 25
$ws.Cells.Item(4,11).Value = 12000
$ws.Cells.Item(4,12).Value = 12000
$ws.Cells.Item(4,13).Value = 12000
$ws.Cells.Item(4,14).Value = "`$/cuna 10 kilos"
$ws.Cells.Item(4,15).Value = "Región Metropolitana"
$ws.Cells.Item(4,16).Value = 1200
$ws.Cells.Item(4,17).Value = 10
$ws.Cells.Item(4,18).Value = "Hortaliza"
